$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: (1,2, p = 15, c = 500) -> (1,2, p = 13, c = 700)
$ws.Range("B2").Value = "(1,2, p = 13, c = 700)"
$ws.Range("C2").Value = 13

# Row 3: (1,4, p = 17, c = 200) -> (1,4, p = 1, c = 133)
$ws.Range("B3").Value = "(1,4, p = 1, c = 133)"
$ws.Range("C3").Value = 1

# Row 4: (2,3, p = 8, c = 100) -> (2,3, p = 10, c = 600)
$ws.Range("B4").Value = "(2,3, p = 10, c = 600)"
$ws.Range("C4").Value = 10

# Row 5: (2,5, p = 7, c = 133) -> (2,5, p = 16, c = 200)
$ws.Range("B5").Value = "(2,5, p = 16, c = 200)"
$ws.Range("C5").Value = 16

# Row 6: (3,6, p = 12, c = 100) -> (3,6, p = 9, c = 200)
$ws.Range("B6").Value = "(3,6, p = 9, c = 200)"
$ws.Range("C6").Value = 9

# Row 7: (4,5, p = 15, c = 700) -> (4,5, p = 10, c = 500)
$ws.Range("B7").Value = "(4,5, p = 10, c = 500)"
$ws.Range("C7").Value = 10

# Row 8: (4,7, p = 1, c = 500) -> (4,7, p = 14, c = 200)
$ws.Range("B8").Value = "(4,7, p = 14, c = 200)"
$ws.Range("C8").Value = 14

# Row 9: (5,6, p = 5, c = 200) -> (5,6, p = 19, c = 700)
$ws.Range("B9").Value = "(5,6, p = 19, c = 700)"
$ws.Range("C9").Value = 19

# Row 10: (5,8, p = 7, c = 500) -> (5,8, p = 4, c = 500)
$ws.Range("B10").Value = "(5,8, p = 4, c = 500)"
$ws.Range("C10").Value = 4

$wb.Save()
